# Add a new worksheet "ODI Batting Extra" after "ODI Batting"
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Header row
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$headerRange = $newSheet.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows
$newSheet.Cells.Item(2, 1).Value = "3544"
$newSheet.Cells.Item(2, 2).Value = 1
$newSheet.Cells.Item(2, 3).Value = "2"
$newSheet.Cells.Item(2, 4).Value = "0"
$newSheet.Cells.Item(2, 5).Value = "8.97%"
$newSheet.Cells.Item(2, 6).Value = "NO"

$newSheet.Cells.Item(3, 1).Value = "3545"
$newSheet.Cells.Item(3, 2).Value = 1
$newSheet.Cells.Item(3, 3).Value = "0"
$newSheet.Cells.Item(3, 4).Value = "0"
$newSheet.Cells.Item(3, 6).Value = "NO"

$newSheet.Cells.Item(4, 1).Value = "3643"
$newSheet.Cells.Item(4, 2).Value = 3
$newSheet.Cells.Item(4, 3).Value = "0"
$newSheet.Cells.Item(4, 4).Value = "0"
$newSheet.Cells.Item(4, 6).Value = "NO"

$newSheet.Cells.Item(5, 1).Value = "3644"
$newSheet.Cells.Item(5, 2).Value = 3
$newSheet.Cells.Item(5, 3).Value = "0"
$newSheet.Cells.Item(5, 4).Value = "0"
$newSheet.Cells.Item(5, 5).Value = "10.48%"
$newSheet.Cells.Item(5, 6).Value = "NO"

$newSheet.Cells.Item(6, 1).Value = "3645"
$newSheet.Cells.Item(6, 6).Value = "NO"
